$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 11:46"

# Update country names that changed position (A column) due to source reordering
$ws.Range("A11").Value = "Suiza"
$ws.Range("A12").Value = "Corea del Sur"
$ws.Range("A14").Value = "Austria"
$ws.Range("A15").Value = "Paises Bajos"
$ws.Range("A51").Value = "Barein"
$ws.Range("A52").Value = "Hong Kong"
$ws.Range("A59").Value = "Libano"
$ws.Range("A60").Value = "Argentina"
$ws.Range("A83").Value = "Albania"
$ws.Range("A84").Value = "Vietnam"
$ws.Range("A85").Value = "Islas Feroe"
$ws.Range("A86").Value = "Republica de Chipre"
$ws.Range("A87").Value = "Burkina Faso"
$ws.Range("A88").Value = "Moldavia"
$ws.Range("A89").Value = "Malta"
$ws.Range("A90").Value = "Brunei"
$ws.Range("A92").Value = "Ucrania"
$ws.Range("A93").Value = "Tunez"
$ws.Range("A94").Value = "Camboya"
$ws.Range("A95").Value = "Senegal"
$ws.Range("A96").Value = "Venezuela"
$ws.Range("A97").Value = "Oman"
$ws.Range("A98").Value = "Bielorrusia"
$ws.Range("A132").Value = "Isla de Man"
$ws.Range("A134").Value = "Guyana"
$ws.Range("A135").Value = "Togo"
$ws.Range("A136").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A137").Value = "Madagascar"
$ws.Range("A138").Value = "Barbados"
$ws.Range("A139").Value = "Kenia"
$ws.Range("A140").Value = "Gibraltar"

# Update numeric statistics (columns B-H) for rows with new data
$ws.Range("B7").Value = 35480
$ws.Range("C7").Value = 344
$ws.Range("E7").Value = 29783
$ws.Range("G7").Value = 31
$ws.Range("H7").Value = 2342

$ws.Range("B8").Value = 30138
$ws.Range("C8").Value = 1082
$ws.Range("E8").Value = 29555

$ws.Range("B11").Value = 9117
$ws.Range("C11").Value = 322
$ws.Range("D11").Value = 131
$ws.Range("E11").Value = 8864
$ws.Range("F11").Value = 141
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 122

$ws.Range("B12").Value = 9037
$ws.Range("C12").Value = 76
$ws.Range("D12").Value = 3507
$ws.Range("E12").Value = 5410
$ws.Range("F12").Value = 59
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 120

$ws.Range("B14").Value = 4767
$ws.Range("C14").Value = 293
$ws.Range("D14").Value = 9
$ws.Range("E14").Value = 4733
$ws.Range("F14").Value = 17
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 25

$ws.Range("B15").Value = 4749
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 4534
$ws.Range("F15").Value = 435
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 213

$ws.Range("B16").Value = 4269
$ws.Range("C16").Value = 526
$ws.Range("E16").Value = 3686

$ws.Range("B17").Value = 2647
$ws.Range("C17").Value = 22
$ws.Range("E17").Value = 2630
$ws.Range("F17").Value = 42

$ws.Range("E27").Value = 1280
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 2

$ws.Range("B31").Value = 908
$ws.Range("C31").Value = 33
$ws.Range("E31").Value = 888
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 7

$ws.Range("F43").Value = 2

$ws.Range("B51").Value = 390
$ws.Range("C51").Value = 13
$ws.Range("D51").Value = 164
$ws.Range("E51").Value = 224
$ws.Range("F51").Value = 3
$ws.Range("H51").Value = 2

$ws.Range("B52").Value = 386
$ws.Range("C52").Value = 29
$ws.Range("D52").Value = 102
$ws.Range("E52").Value = 280
$ws.Range("F52").Value = 4
$ws.Range("H52").Value = 4

$ws.Range("D53").Value = 7
$ws.Range("E53").Value = 362

$ws.Range("B59").Value = 304
$ws.Range("C59").Value = 37
$ws.Range("D59").Value = 8
$ws.Range("E59").Value = 292
$ws.Range("F59").Value = 4

$ws.Range("B60").Value = 301
$ws.Range("D60").Value = 51
$ws.Range("E60").Value = 246
$ws.Range("F60").Value = 0

$ws.Range("E74").Value = 157
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 9

$ws.Range("C83").Value = 19
$ws.Range("D83").Value = 10
$ws.Range("E83").Value = 108
$ws.Range("F83").Value = 2
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 5

$ws.Range("B84").Value = 123
$ws.Range("D84").Value = 17
$ws.Range("E84").Value = 106
$ws.Range("F84").Value = 3

$ws.Range("B85").Value = 122
$ws.Range("C85").Value = 4
$ws.Range("D85").Value = 23
$ws.Range("E85").Value = 99
$ws.Range("F85").Value = 0
$ws.Range("H85").Value = 0

$ws.Range("B86").Value = 116
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 3
$ws.Range("E86").Value = 112
$ws.Range("F86").Value = 3
$ws.Range("H86").Value = 1

$ws.Range("B87").Value = 114
$ws.Range("C87").Value = 15
$ws.Range("D87").Value = 7
$ws.Range("E87").Value = 103
$ws.Range("F87").Value = 0
$ws.Range("H87").Value = 4

$ws.Range("B88").Value = 109
$ws.Range("E88").Value = 106
$ws.Range("F88").Value = 10
$ws.Range("H88").Value = 1

$ws.Range("B89").Value = 107
$ws.Range("C89").Value = 0
$ws.Range("E89").Value = 105
$ws.Range("F89").Value = 1

$ws.Range("C90").Value = 13
$ws.Range("E90").Value = 102
$ws.Range("H90").Value = 0

$ws.Range("B92").Value = 97
$ws.Range("C92").Value = 24
$ws.Range("E92").Value = 93
$ws.Range("F92").Value = 0

$ws.Range("B93").Value = 90
$ws.Range("C93").Value = 1
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = 86
$ws.Range("F93").Value = 11
$ws.Range("H93").Value = 3

$ws.Range("B94").Value = 87
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 2
$ws.Range("E94").Value = 85
$ws.Range("F94").Value = 1
$ws.Range("H94").Value = 0

$ws.Range("B95").Value = 86
$ws.Range("C95").Value = 7
$ws.Range("D95").Value = 8
$ws.Range("E95").Value = 78
$ws.Range("F95").Value = 0

$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 15
$ws.Range("E96").Value = 69
$ws.Range("F96").Value = 2

$ws.Range("B97").Value = 84
$ws.Range("C97").Value = 18
$ws.Range("D97").Value = 17
$ws.Range("E97").Value = 67

$ws.Range("B98").Value = 81
$ws.Range("D98").Value = 22
$ws.Range("E98").Value = 59

$ws.Range("C132").Value = 7
$ws.Range("E132").Value = 20
$ws.Range("H132").Value = 0

$ws.Range("B134").Value = 20
$ws.Range("E134").Value = 19
$ws.Range("H134").Value = 1

$ws.Range("B135").Value = 18
$ws.Range("E135").Value = 18

$ws.Range("C136").Value = 0

$ws.Range("C137").Value = 5

$ws.Range("B138").Value = 17
$ws.Range("E138").Value = 17

$ws.Range("B139").Value = 16
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 16

$ws.Range("B140").Value = 15
$ws.Range("D140").Value = 5
$ws.Range("E140").Value = 10

